# Update format for output:
#  - Rows 2-15 and 20 get their "shock" column reformatted from raw decimal
#    numbers into percentage-style text labels (and some rows gain/lose an
#    "extreme_level" label too).
#  - Old rows 16-20 (M15..M19) are dropped, and the tail of the table
#    (previously M20..M30, rows 21-31) is replaced by a shorter, relabeled
#    block M20..M25 occupying rows 16-21.
#  - The sheet shrinks from A1:D31 down to A1:D21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($range, $text)
    # Force a literal-text number format first so values that look like
    # numbers/percentages (e.g. "0.1 %") are not silently reinterpreted
    # by Excel as numeric percentages.
    $range.NumberFormat = "@"
    $range.Value = $text
}

# --- Remove the trailing rows that disappear entirely (old rows 22-31) ---
$ws.Rows("22:31").Delete()

# --- Reformat the "shock" values for rows 2-10 (M1..M9) ---
Set-TextCell $ws.Range("C2")  "0.1 %"
Set-TextCell $ws.Range("C3")  "0.1 %"
Set-TextCell $ws.Range("C4")  "0.1 %"
Set-TextCell $ws.Range("C5")  "0.1 %"
Set-TextCell $ws.Range("C6")  "0.1 %"
Set-TextCell $ws.Range("C7")  "0.1 %"
Set-TextCell $ws.Range("C8")  "0.1 %"
Set-TextCell $ws.Range("C9")  "0.1 %"
Set-TextCell $ws.Range("C10") "0.1 %"

# Row 11 (M10)
Set-TextCell $ws.Range("C11") "0.8 %"

# Rows 12-15 (M11..M14) gain both a "shock" label and an "extreme_level" label
Set-TextCell $ws.Range("C12") "13% max"
Set-TextCell $ws.Range("D12") "(+1746 ppts)"

Set-TextCell $ws.Range("C13") "13 peak"
Set-TextCell $ws.Range("D13") "(+1764 ppts)"

Set-TextCell $ws.Range("C14") "13% peak"
Set-TextCell $ws.Range("D14") "(+178200 bps)"

Set-TextCell $ws.Range("C15") "13% peak"
Set-TextCell $ws.Range("D15") "(+180000 bps)"

# --- Replace rows 16-21 (formerly the M15-M19 block plus the old M20 row)
#     with the relabeled M20-M25 block ---
Set-TextCell $ws.Range("A16") "M20"
Set-TextCell $ws.Range("B16") "MMM20"
Set-TextCell $ws.Range("C16") "0.2 %"
Set-TextCell $ws.Range("D16") "0.2 %"

Set-TextCell $ws.Range("A17") "M21"
Set-TextCell $ws.Range("B17") "MMM21"
Set-TextCell $ws.Range("C17") "1% trough"
Set-TextCell $ws.Range("D17") "(191400 bps)"

Set-TextCell $ws.Range("A18") "M22"
Set-TextCell $ws.Range("B18") "MMM22"
Set-TextCell $ws.Range("C18") "1% trough"
Set-TextCell $ws.Range("D18") "(193200 bps)"

Set-TextCell $ws.Range("A19") "M23"
Set-TextCell $ws.Range("B19") "MMM23"
Set-TextCell $ws.Range("C19") "1% trough"
Set-TextCell $ws.Range("D19") "(195000 bps)"

Set-TextCell $ws.Range("A20") "M24"
Set-TextCell $ws.Range("B20") "MMM24"
Set-TextCell $ws.Range("C20") "0.1 %"
Set-TextCell $ws.Range("D20") ""

Set-TextCell $ws.Range("A21") "M25"
Set-TextCell $ws.Range("B21") "MMM25"
Set-TextCell $ws.Range("C21") "1% trough"
Set-TextCell $ws.Range("D21") "(198600 bps)"
